$wb = $excel.ActiveWorkbook

# --- Update config_netNodes sheet (re-order execution of MS/HS/LS nodes) ---
# EnergyBalances must be calculated 'bottom up' on the grid, so the LS nodes
# (E2, E3, E4) need to execute before their MS parent (E1), which in turn
# needs to execute before the HS node (H1).
$ws = $wb.Worksheets.Item("config_netNodes")

# Row 2 -> data for E2 (previously in row 3)
$ws.Range("C2").Value = "E2"
$ws.Range("D2").Value = "ELECTRICITY"
$ws.Range("E2").Value = "MSLS"
$ws.Range("F2").Value = "E1"
$ws.Range("G2").Value = 1000

# Row 3 -> data for E3 (previously in row 4)
$ws.Range("C3").Value = "E3"
$ws.Range("D3").Value = "ELECTRICITY"
$ws.Range("E3").Value = "MSLS"
$ws.Range("F3").Value = "E1"
$ws.Range("G3").Value = 1000

# Row 4 -> data for E4 (previously in row 6)
$ws.Range("C4").Value = "E4"
$ws.Range("D4").Value = "ELECTRICITY"
$ws.Range("E4").Value = "MSLS"
$ws.Range("F4").Value = "E1"
$ws.Range("G4").Value = 1000

# Row 5 -> data for E1 (previously in row 2); no parent column value anymore
$ws.Range("C5").Value = "E1"
$ws.Range("D5").Value = "ELECTRICITY"
$ws.Range("E5").Value = "HSMS"
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = 500000

# Row 6 -> data for H1 (previously in row 5); no parent column value anymore
$ws.Range("C6").Value = "H1"
$ws.Range("D6").Value = "HEAT"
$ws.Range("E6").Value = "MT"
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = 300000

# --- Sheet view / selection adjustments ---
# config_actors keeps its own last selection, but it is no longer the
# active sheet, so it loses the tabSelected flag once another sheet is
# activated below.
$actors = $wb.Worksheets.Item("config_actors")
$actors.Activate()
$actors.Range("E22").Select()

# config_netNodes becomes the active/selected sheet with a new selection.
$ws.Activate()
$ws.Range("D11").Select()
